$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CHECKOUT_04 (Giỏ hàng rỗng / empty cart) ---
$ws.Range("A2").Value = "CHECKOUT_04"
$ws.Range("B2").Value = "Giỏ hàng rỗng"
$ws.Range("C2").Value = "Cart=[]"
$ws.Range("D2").Value = "1. Service return EMPTY_CART"
$ws.Range("E2").Value = "Redirect view-cart.jsp"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# --- Row 3: CHECKOUT_01 (Đặt hàng thành công / success) ---
$ws.Range("A3").Value = "CHECKOUT_01"
$ws.Range("B3").Value = "Đặt hàng thành công"
$ws.Range("C3").Value = "Addr: Hanoi"
$ws.Range("D3").Value = "1. Service return SUCCESS"
$ws.Range("E3").Value = "Xóa Cart & Redirect Home"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# --- Row 4: CHECKOUT_05 (Lỗi lưu Order (DB) / order save error) ---
$ws.Range("A4").Value = "CHECKOUT_05"
$ws.Range("B4").Value = "Lỗi lưu Order (DB)"
$ws.Range("C4").Value = "DB Error"
$ws.Range("D4").Value = "1. Service return ORDER_FAILED"
$ws.Range("E4").Value = "Forward Checkout.jsp + Error"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

# --- Row 5 (new): CHECKOUT_03 (Thiếu thông tin / missing info) ---
$ws.Range("A5").Value = "CHECKOUT_03"
$ws.Range("B5").Value = "Thiếu thông tin"
$ws.Range("C5").Value = "Addr=null"
$ws.Range("D5").Value = "1. Service return MISSING_INFO"
$ws.Range("E5").Value = "Forward Checkout.jsp + Error"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# --- Row 6 (new): CHECKOUT_02 (Chưa đăng nhập / not logged in) ---
$ws.Range("A6").Value = "CHECKOUT_02"
$ws.Range("B6").Value = "Chưa đăng nhập"
$ws.Range("C6").Value = "User=null"
$ws.Range("D6").Value = "1. Service return NOT_LOGGED_IN"
$ws.Range("E6").Value = "Redirect Login"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

# --- Row 7 (new): CHECKOUT_07 (Lỗi hệ thống / system error) ---
$ws.Range("A7").Value = "CHECKOUT_07"
$ws.Range("B7").Value = "Lỗi hệ thống"
$ws.Range("C7").Value = "Crash"
$ws.Range("D7").Value = "1. Service return EXCEPTION"
$ws.Range("E7").Value = "Forward Checkout.jsp + Error"
$ws.Range("F7").Value = "OK"
$ws.Range("G7").Value = "PASS"

# --- Row 8 (new): CHECKOUT_06 (Lỗi lưu chi tiết / detail save error) ---
$ws.Range("A8").Value = "CHECKOUT_06"
$ws.Range("B8").Value = "Lỗi lưu chi tiết"
$ws.Range("C8").Value = "DB Error"
$ws.Range("D8").Value = "1. Service return DETAIL_FAILED"
$ws.Range("E8").Value = "Forward Checkout.jsp + Error"
$ws.Range("F8").Value = "OK"
$ws.Range("G8").Value = "PASS"

# New rows 5-8 need the same "PASS" cell formatting (bold green) that column G
# already carries on rows 2-4 -- copy formats only so the style table isn't
# duplicated with new near-identical entries.
$ws.Range("G4").Copy()
$ws.Range("G5:G8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column C/D/E need to be narrower/wider to fit the new data (bestFit).
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 31.166666666666668
$ws.Columns.Item(5).ColumnWidth = 26.333333333333332
